$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# Widen column B so the new, longer row label fits
$ws.Columns.Item(2).ColumnWidth = 23

# Clear the stray "2025" value that was sitting in Q3
$ws.Range("Q3").Value = $null

# Row 7 "Company" - quarterly company-owned shop revenue ($mm)
$ws.Range("J7").Value = 295
$ws.Range("K7").Value = 308
$ws.Range("N7").Value = 381
$ws.Range("O7").Value = 393

# Row 8 "Franchise" - quarterly franchise revenue ($mm)
$ws.Range("J8").Value = 30
$ws.Range("K8").Value = 30
$ws.Range("N8").Value = 35
$ws.Range("O8").Value = 31

# Row 5 "Revenue/Shop" = Revenue / Shops
$ws.Range("C5:J5").Formula = "=C9/C4"
$ws.Range("K5").Formula = "=K9/K4"
$ws.Range("L5:O5").Formula = "=L9/L4"

# Row 9 "Revenue" = Company + Franchise, bold
$ws.Range("C9").Formula = "=C7+C8"
$ws.Range("D9:N9").Formula = "=D7+D8"
$ws.Range("O9").Formula = "=O7+O8"
$ws.Range("C9:O9").Font.Bold = $true

# Number format for the Revenue/Shop growth row
$ws.Range("C5:O5").NumberFormat = "0.000"

# New label row for the revenue growth metric
$ws.Range("B21").Value = "Revenue Growth y/o/y %"
$ws.Range("B21").Font.ThemeColor = 1
$ws.Rows.Item(21).RowHeight = 15.75

# Restore the active-cell selection on the frozen pane
$ws.Range("C4").Select()
